$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.925.98'
$ws.Range("E2").Value = '  -0.98%  '
$ws.Range("D3").Value = '3.388.78'
$ws.Range("E3").Value = '  -1.51%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '572.05'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.88%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '142.15'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.68%  '
$ws.Range("E7").Value = '  +0.05%  '
$ws.Range("D8").Value = '3.388.97'
$ws.Range("E8").Value = '  -1.52%  '
$ws.Range("E9").Value = '  -0.36%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.54'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.95%  '
$ws.Range("E11").Value = '  -2.47%  '
$ws.Range("E12").Value = '  +2.02%  '
$ws.Range("D13").Value = '3.965.93'
$ws.Range("E13").Value = '  -1.48%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '28.19'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.78%  '
$ws.Range("E16").Value = '  -1.53%  '
$ws.Range("D17").Value = '3.388.97'
$ws.Range("E17").Value = '  -1.55%  '
$ws.Range("D18").Value = '60.972.41'
$ws.Range("E18").Value = '  -1.03%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.17'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.58%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.82'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -3.33%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '8.97'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -4.95%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '384.13'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -3.63%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.558'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.76%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '74.36'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.74%  '
$ws.Range("E25").Value = '  +0.36%  '
$ws.Range("E26").Value = '  -4.99%  '
$ws.Range("D27").Value = '3.524.31'
$ws.Range("E27").Value = '  -1.46%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.180'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.77%  '
$ws.Range("E29").Value = '  -0.09%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.39'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.09%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.98'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -3.27%  '
$ws.Range("E32").Value = '  -1.96%  '
$ws.Range("E33").Value = '  -2.51%  '
$ws.Range("E34").Value = '  +0.04%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '23.52'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -2.00%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.99'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.69%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '167.77'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.33%  '
$ws.Range("D38").Value = '3.416.77'
$ws.Range("E38").Value = '  -1.43%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '4.99'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -3.03%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.49'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -4.55%  '
$ws.Range("B41").Value = 'Hedera'
$ws.Range("C41").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0775'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.97%  '
$ws.Range("B42").Value = 'EnergySwap'
$ws.Range("C42").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '27.60'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.07%  '
$ws.Range("E43").Value = '  +0.04%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.780'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.76%  '
$ws.Range("E45").Value = '  -0.39%  '
$ws.Range("E46").Value = '  -1.54%  '
$ws.Range("E47").Value = '  -3.78%  '
$ws.Range("E48").Value = '  -0.42%  '
$ws.Range("D49").Value = '2.478.69'
$ws.Range("E49").Value = '  -4.85%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '6.82'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.83%  '
$ws.Range("E51").Value = '  -0.87%  '
